# Refresh crypto price/volume snapshot (GitHub Actions scheduled update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "'39.741.12"
$ws.Range("E2").Value = "  +2.56%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "'2.162.01"
$ws.Range("E3").Value = "  +2.95%  "

# Row 4: TetherUSD
$ws.Range("E4").Value = "  -0.03%  "

# Row 5: BNB
$ws.Range("D5").Value = "'226.79"
$ws.Range("E5").Value = "  -0.15%  "

# Row 6: XRP
$ws.Range("E6").Value = "  +1.10%  "

# Row 7: Solana
$ws.Range("D7").Value = "'62.91"
$ws.Range("E7").Value = "  +1.76%  "

# Row 8: USDC
$ws.Range("E8").Value = "  +0.03%  "

# Row 9: Cardano
$ws.Range("E9").Value = "  +0.85%  "

# Row 10: Dogecoin
$ws.Range("D10").Value = "'0.0842"
$ws.Range("E10").Value = "  +0.18%  "

# Row 11: TRON
$ws.Range("E11").Value = "  +0.34%  "

# Row 12: Chainlink
$ws.Range("E12").Value = "  +0.82%  "

# Row 13: WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "'2.483.02"
$ws.Range("E13").Value = "  +3.03%  "

# Row 14: Avalanche
$ws.Range("D14").Value = "'21.78"
$ws.Range("E14").Value = "  -0.78%  "

# Row 15: Polygon
$ws.Range("E15").Value = "  +0.97%  "

# Row 16: Polkadot
$ws.Range("D16").Value = "'5.50"
$ws.Range("E16").Value = "  +0.50%  "

# Row 17: WrappedEther
$ws.Range("D17").Value = "'2.158.21"
$ws.Range("E17").Value = "  +2.42%  "

# Row 18: WrappedBTC
$ws.Range("D18").Value = "'39.749.50"
$ws.Range("E18").Value = "  +2.73%  "

# Row 19: Litecoin
$ws.Range("D19").Value = "'71.79"
$ws.Range("E19").Value = "  +0.30%  "

# Row 20: Uniswap
$ws.Range("E20").Value = "  -0.19%  "

# Row 21: ShibaInu
$ws.Range("D21").Value = "0.0₃0848"
$ws.Range("E21").Value = "  +0.79%  "

# Row 22: BitcoinCash
$ws.Range("D22").Value = "'228.21"
$ws.Range("E22").Value = "  +0.65%  "

# Row 23: Dai
$ws.Range("E23").Value = "  +0.10%  "

# Row 24: PancakeSwap
$ws.Range("E24").Value = "  +2.20%  "

# Row 25: Toncoin
$ws.Range("D25").Value = "'2.33"
$ws.Range("E25").Value = "  -1.20%  "

# Row 26: Monero
$ws.Range("D26").Value = "'172.05"
$ws.Range("E26").Value = "  +1.07%  "

# Row 27: Cosmos
$ws.Range("E27").Value = "  -2.11%  "

# Row 28: Kaspa
$ws.Range("E28").Value = "  +2.28%  "

# Row 29: ImmutableX
$ws.Range("E29").Value = "  +1.70%  "

# Row 30: EthereumClassic
$ws.Range("D30").Value = "'19.61"
$ws.Range("E30").Value = "  +1.44%  "

# Row 31: WEMIXToken
$ws.Range("E31").Value = "  +6.06%  "

# Row 32: Stellar
$ws.Range("E32").Value = "  +0.83%  "

# Row 33: Filecoin
$ws.Range("E33").Value = "  +0.69%  "

# Row 34: InternetComputer(DFINITY)
$ws.Range("E34").Value = "  -2.52%  "

# Row 35: THORChain
$ws.Range("D35").Value = "'6.94"
$ws.Range("E35").Value = "  -2.48%  "

# Row 36: Hedera
$ws.Range("D36").Value = "'0.0617"
$ws.Range("E36").Value = "  +0.68%  "

# Row 37: RenderToken
$ws.Range("D37").Value = "'3.76"
$ws.Range("E37").Value = "  +7.93%  "

# Row 38: LidoDAOToken
$ws.Range("D38").Value = "'2.40"
$ws.Range("E38").Value = "  +1.91%  "

# Row 39: BinanceUSD
$ws.Range("D39").Value = "'0.999"
$ws.Range("E39").Value = "  -0.13%  "

# Row 40: FTXToken
$ws.Range("D40").Value = "'4.92"
$ws.Range("E40").Value = "  +17.74%  "

# Row 41: Aave
$ws.Range("D41").Value = "'102.85"
$ws.Range("E41").Value = "  +1.28%  "

# Row 42: VeChain
$ws.Range("E42").Value = "  -0.73%  "

# Row 43: InjectiveProtocol
$ws.Range("E43").Value = "  -2.28%  "

# Row 44: Maker
$ws.Range("D44").Value = "'1.511.86"
$ws.Range("E44").Value = "  -0.83%  "

# Row 45: TrustWalletToken
$ws.Range("E45").Value = "  -0.15%  "

# Row 46: FraxShare
$ws.Range("E46").Value = "  +2.33%  "

# Row 47: Cronos
$ws.Range("D47").Value = "'0.0930"
$ws.Range("E47").Value = "  +2.30%  "

# Row 48: HuobiToken
$ws.Range("D48").Value = "'2.80"
$ws.Range("E48").Value = "  -0.02%  "

# Row 49: ARBITRUM
$ws.Range("E49").Value = "  +1.17%  "

# Row 50: ranking swap, MultiversX -> TerraClassic
$ws.Range("B50").Value = "TerraClassic"
$ws.Range("C50").Value = "https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc"
$ws.Range("D50").Value = "'0.000192"
$ws.Range("E50").Value = "  +28.88%  "

# Row 51: ranking swap, TerraClassic -> MultiversX
$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D51").Value = "'49.68"
$ws.Range("E51").Value = "  +8.39%  "
